$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.63
$ws.Range("G2").Value = 1.84
$ws.Range("H2").Value = 5.3
$ws.Range("I2").Value = 6.8
$ws.Range("K2").Value = 4.3
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 1.85
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 2.18
$ws.Range("Z2").Value = 120
$ws.Range("AB2").Value = 9
$ws.Range("AL2").Value = 110
$ws.Range("AN2").Value = 12.5
$ws.Range("F3").Value = 1.58
$ws.Range("G3").Value = 1.59
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 7.2
$ws.Range("J3").Value = 4.3
$ws.Range("L3").Value = 1.33
$ws.Range("N3").Value = 4.3
$ws.Range("O3").Value = 1.24
$ws.Range("P3").Value = 2.2
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.46
$ws.Range("S3").Value = 2.78
$ws.Range("U3").Value = 2.08
$ws.Range("W3").Value = 2.7
$ws.Range("X3").Value = 19
$ws.Range("Z3").Value = 130
$ws.Range("AA3").Value = 190
$ws.Range("AC3").Value = 11
$ws.Range("AF3").Value = 10
$ws.Range("AG3").Value = 9.800000000000001
$ws.Range("AI3").Value = 210
$ws.Range("AK3").Value = 16
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 580
$ws.Range("AN3").Value = 8
$ws.Range("AO3").Value = 100
$ws.Range("F4").Value = 2
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 3.7
$ws.Range("V4").Value = 1.33
$ws.Range("X4").Value = 20
$ws.Range("Y4").Value = 17.5
$ws.Range("Z4").Value = 32
$ws.Range("AA4").Value = 190
$ws.Range("AB4").Value = 12
$ws.Range("AC4").Value = 9.4
$ws.Range("AD4").Value = 16.5
$ws.Range("AF4").Value = 15
$ws.Range("F5").Value = 2.26
$ws.Range("H5").Value = 3.5
$ws.Range("N5").Value = 2.76
$ws.Range("V5").Value = 1.33
$ws.Range("W5").Value = 1.69
$ws.Range("X5").Value = 17.5
$ws.Range("Y5").Value = 12
$ws.Range("Z5").Value = 26
$ws.Range("AC5").Value = 7.6
$ws.Range("AD5").Value = 32
$ws.Range("AG5").Value = 12
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 150
$ws.Range("AO5").Value = 90
$ws.Range("F6").Value = 2.48
$ws.Range("K6").Value = 3.3
$ws.Range("X6").Value = 18.5
$ws.Range("AC6").Value = 17.5
$ws.Range("F7").Value = 3.25
$ws.Range("G7").Value = 4.7
$ws.Range("H7").Value = 2.14
$ws.Range("K7").Value = 3.55
$ws.Range("L7").Value = 1.43
$ws.Range("V7").Value = 1.7
$ws.Range("W7").Value = 1.3
